$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1551.5
$ws.Range("I43").Value = 1639.5
$ws.Range("J43").Value = 1199.5
$ws.Range("K43").Value = 1639.5
$ws.Range("L43").Value = 1199.5
$ws.Range("M43").Value = -1570.5
$ws.Range("N43").Value = -1337.5

$ws.Range("H74").Value = 10050.1
$ws.Range("I74").Value = 10050.1
$ws.Range("K74").Value = 10050.1
$ws.Range("M74").Value = -9114.1

$ws.Range("H76").Value = 4031.5625
$ws.Range("I76").Value = 3923.3845
$ws.Range("K76").Value = 3923.3845
$ws.Range("M76").Value = -3608.3845

$ws.Range("H77").Value = 10050.1
$ws.Range("I77").Value = 10050.1
$ws.Range("K77").Value = 50250.5
$ws.Range("M77").Value = -45570.5

$ws.Range("H79").Value = 4031.5625
$ws.Range("I79").Value = 3923.3845
$ws.Range("K79").Value = 3923.3845
$ws.Range("M79").Value = -2831.3845

$ws.Range("H86").Value = 5216.0835
$ws.Range("I86").Value = 5666.5
$ws.Range("J86").Value = 4765.6665
$ws.Range("K86").Value = 5666.5
$ws.Range("L86").Value = 4765.6665
$ws.Range("M86").Value = -4543.5
$ws.Range("N86").Value = -7011.6665

$ws.Range("H89").Value = 5216.0835
$ws.Range("I89").Value = 5666.5
$ws.Range("J89").Value = 4765.6665
$ws.Range("K89").Value = 28332.5
$ws.Range("L89").Value = 23828.3325
$ws.Range("M89").Value = -22716.5
$ws.Range("N89").Value = -35060.3325

$ws.Range("H103").Value = 1708.3334
$ws.Range("I103").Value = 1750
$ws.Range("K103").Value = 5250
$ws.Range("M103").Value = -4664

$ws.Range("H132").Value = 2679.4443
$ws.Range("I132").Value = 2834.375
$ws.Range("K132").Value = 8503.125
$ws.Range("M132").Value = -5973.125

$ws.Range("H137").Value = 10008651
$ws.Range("I137").Value = 18184598
$ws.Range("J137").Value = 15828.444
$ws.Range("K137").Value = 54553794
$ws.Range("L137").Value = 47485.33199999999
$ws.Range("M137").Value = -54551244
$ws.Range("N137").Value = -52585.33199999999

$ws.Range("H138").Value = 6159.936
$ws.Range("I138").Value = 3972
$ws.Range("J138").Value = 7185.5312
$ws.Range("K138").Value = 11916
$ws.Range("L138").Value = 21556.5936
$ws.Range("M138").Value = -6776
$ws.Range("N138").Value = -31836.5936

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1044307.5
$ws.Range("I32").Value = 1283748
$ws.Range("K32").Value = 1283748
$ws.Range("M32").Value = -1283461

$ws.Range("H97").Value = 2004.25
$ws.Range("I97").Value = 2004.25
$ws.Range("K97").Value = 2004.25
$ws.Range("M97").Value = -1508.25

$ws.Range("M110").ClearContents()
$ws.Range("H110").Value = 1990
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 1990
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 1990
$ws.Range("N110").Value = -6080

$ws.Range("N125").ClearContents()
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0

$ws.Range("H132").Value = 5272.907
$ws.Range("I132").Value = 3037.64
$ws.Range("K132").Value = 9112.92
$ws.Range("M132").Value = -6582.92

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N76").ClearContents()
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0

$ws.Range("N79").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0

$ws.Range("H94").Value = 2835.6155
$ws.Range("I94").Value = 3077.6843
$ws.Range("J94").Value = 2178.5715
$ws.Range("K94").Value = 3077.6843
$ws.Range("L94").Value = 2178.5715
$ws.Range("M94").Value = -2626.6843
$ws.Range("N94").Value = -3080.5715

$ws.Range("H134").Value = 12596254
$ws.Range("J134").Value = 18527352
$ws.Range("L134").Value = 55582056
$ws.Range("N134").Value = -55587126

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 402235.75
$ws.Range("I31").Value = 824489.5
$ws.Range("J31").Value = 4820.4414
$ws.Range("K31").Value = 824489.5
$ws.Range("L31").Value = 4820.4414
$ws.Range("M31").Value = -824194.5
$ws.Range("N31").Value = -5410.4414

$ws.Range("H34").Value = 402235.75
$ws.Range("I34").Value = 824489.5
$ws.Range("J34").Value = 4820.4414
$ws.Range("K34").Value = 824489.5
$ws.Range("L34").Value = 4820.4414
$ws.Range("M34").Value = -824287.5
$ws.Range("N34").Value = -5224.4414

$ws.Range("H58").Value = 19108358
$ws.Range("I58").Value = 27782696
$ws.Range("K58").Value = 27782696
$ws.Range("M58").Value = -27782493

$ws.Range("H132").Value = 5738.25
$ws.Range("I132").Value = 5704.9165
$ws.Range("K132").Value = 17114.7495
$ws.Range("M132").Value = -14584.7495

$ws.Range("H136").Value = 19108358
$ws.Range("I136").Value = 27782696
$ws.Range("K136").Value = 83348088
$ws.Range("M136").Value = -83345538

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2184632
$ws.Range("I5").Value = 1786197.2
$ws.Range("K5").Value = 5358591.6
$ws.Range("M5").Value = -5358479.6

$ws.Range("H7").Value = 650.25
$ws.Range("I7").Value = 200.33333
$ws.Range("K7").Value = 600.99999
$ws.Range("M7").Value = -488.99999

$ws.Range("H10").Value = 1673.3334
$ws.Range("I10").Value = 10
$ws.Range("K10").Value = 30
$ws.Range("M10").Value = 109

$ws.Range("H12").Value = 849.4
$ws.Range("J12").Value = 641.25
$ws.Range("L12").Value = 1923.75
$ws.Range("N12").Value = -2269.75

$ws.Range("H113").Value = 385.25
$ws.Range("J113").Value = 362.72726
$ws.Range("L113").Value = 1088.18178
$ws.Range("N113").Value = -5428.18178

$ws.Range("H135").Value = 2184632
$ws.Range("I135").Value = 1786197.2
$ws.Range("K135").Value = 16075774.8
$ws.Range("M135").Value = -16073239.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 1126.6666
$ws.Range("I32").Value = 4833.3335
$ws.Range("J32").Value = 200
$ws.Range("K32").Value = 4833.3335
$ws.Range("L32").Value = 200
$ws.Range("M32").Value = -4516.3335
$ws.Range("N32").Value = -834

$ws.Range("H40").Value = 5374.5625
$ws.Range("I40").Value = 5199.533
$ws.Range("J40").Value = 8000
$ws.Range("K40").Value = 5199.533
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = -5063.533
$ws.Range("N40").Value = -8272

$ws.Range("H61").Value = 7112.522
$ws.Range("I61").Value = 7896
$ws.Range("K61").Value = 7896
$ws.Range("M61").Value = -7694

$ws.Range("H68").Value = 4998.5
$ws.Range("I68").Value = 4998.5
$ws.Range("K68").Value = 4998.5
$ws.Range("M68").Value = -4249.5

$ws.Range("H71").Value = 4998.5
$ws.Range("I71").Value = 4998.5
$ws.Range("K71").Value = 24992.5
$ws.Range("M71").Value = -21248.5

$ws.Range("H113").Value = 7112.522
$ws.Range("I113").Value = 7896
$ws.Range("K113").Value = 7896
$ws.Range("M113").Value = -5726

$ws.Range("H122").Value = 6591.75
$ws.Range("I122").Value = 6463.727
$ws.Range("K122").Value = 19391.181
$ws.Range("M122").Value = -16941.181

$ws.Range("H132").Value = 2384537.2
$ws.Range("I132").Value = 3707324.5
$ws.Range("J132").Value = 3520
$ws.Range("K132").Value = 11121973.5
$ws.Range("L132").Value = 10560
$ws.Range("M132").Value = -11119443.5
$ws.Range("N132").Value = -15620

$ws.Range("H136").Value = 19178558
$ws.Range("I136").Value = 13902101
$ws.Range("J136").Value = 66666664
$ws.Range("K136").Value = 41706303
$ws.Range("L136").Value = 199999992
$ws.Range("M136").Value = -41703753
$ws.Range("N136").Value = -200005092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M11").ClearContents()
$ws.Range("H11").Value = 50000000
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0

$ws.Range("H107").Value = 454.84
$ws.Range("I107").Value = 485.8421
$ws.Range("J107").Value = 356.66666
$ws.Range("K107").Value = 1457.5263
$ws.Range("L107").Value = 1069.99998
$ws.Range("M107").Value = 462.4737
$ws.Range("N107").Value = -4909.999980000001

$ws.Range("H122").Value = 65459.2
$ws.Range("I122").Value = 7866.75
$ws.Range("K122").Value = 23600.25
$ws.Range("M122").Value = -21150.25

$ws.Range("H132").Value = 5378718.5
$ws.Range("I132").Value = 6946817.5
$ws.Range("J132").Value = 2379
$ws.Range("K132").Value = 20840452.5
$ws.Range("L132").Value = 7137
$ws.Range("M132").Value = -20837922.5
$ws.Range("N132").Value = -12197
